$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71
$ws.Range("B71").Value = 6139072
$ws.Range("E71").Value = "JK Tammeka Tartu"
$ws.Range("F71").Value = "FC Flora Tallinn"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 2
$ws.Range("I71").Value = "A"
$ws.Range("J71").Value = 9
$ws.Range("K71").Value = 7
$ws.Range("L71").Value = 1.166
$ws.Range("M71").Value = 7
$ws.Range("N71").Value = 6
$ws.Range("O71").Value = 1.25
$ws.Range("P71").Value = 1.75
$ws.Range("Q71").Value = 1.9
$ws.Range("R71").Value = 1.9
$ws.Range("S71").Value = 3
$ws.Range("T71").Value = 1.95
$ws.Range("U71").Value = 1.85
$ws.Range("V71").Value = -1
$ws.Range("W71").Value = -1
$ws.Range("X71").Value = 0.25
$ws.Range("Y71").Value = 0.8999999999999999
$ws.Range("Z71").Value = -1
$ws.Range("AA71").Value = 0
$ws.Range("AB71").Value = 0

# Row 72
$ws.Range("B72").Value = 6139071
$ws.Range("E72").Value = "Parnu JK Vaprus"
$ws.Range("F72").Value = "JK Trans Narva"
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 2
$ws.Range("I72").Value = "H"
$ws.Range("J72").Value = 2.4
$ws.Range("K72").Value = 3.2
$ws.Range("L72").Value = 2.6
$ws.Range("M72").Value = 3
$ws.Range("N72").Value = 3.25
$ws.Range("O72").Value = 2.2
$ws.Range("P72").Value = 0.25
$ws.Range("Q72").Value = 1.825
$ws.Range("R72").Value = 1.975
$ws.Range("S72").Value = 2.5
$ws.Range("T72").Value = 1.875
$ws.Range("U72").Value = 1.925
$ws.Range("V72").Value = 2
$ws.Range("W72").Value = -1
$ws.Range("X72").Value = -1
$ws.Range("Y72").Value = 0.825
$ws.Range("Z72").Value = -1
$ws.Range("AA72").Value = 0.875
$ws.Range("AB72").Value = -1

# Row 95
$ws.Range("B95").Value = 6416370
$ws.Range("E95").Value = "FC Levadia Tallinn"
$ws.Range("F95").Value = "Parnu JK Vaprus"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = "D"
$ws.Range("J95").Value = 1.166
$ws.Range("K95").Value = 7
$ws.Range("L95").Value = 11
$ws.Range("M95").Value = 1.2
$ws.Range("N95").Value = 6
$ws.Range("O95").Value = 11
$ws.Range("P95").Value = -2
$ws.Range("Q95").Value = 1.85
$ws.Range("R95").Value = 1.95
$ws.Range("S95").Value = 3
$ws.Range("T95").Value = 1.85
$ws.Range("U95").Value = 1.95
$ws.Range("V95").Value = -1
$ws.Range("W95").Value = 5
$ws.Range("X95").Value = -1
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.95
$ws.Range("AA95").Value = -1
$ws.Range("AB95").Value = 0.95

# Row 96
$ws.Range("B96").Value = 6482819
$ws.Range("E96").Value = "JK Tammeka Tartu"
$ws.Range("F96").Value = "FC Kuressaare"
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 1
$ws.Range("I96").Value = "A"
$ws.Range("J96").Value = 1.833
$ws.Range("K96").Value = 3.5
$ws.Range("L96").Value = 3.5
$ws.Range("M96").Value = 2.1
$ws.Range("N96").Value = 3.4
$ws.Range("O96").Value = 2.875
$ws.Range("P96").Value = -0.25
$ws.Range("Q96").Value = 1.975
$ws.Range("R96").Value = 1.825
$ws.Range("S96").Value = 3
$ws.Range("T96").Value = 1.825
$ws.Range("U96").Value = 1.975
$ws.Range("V96").Value = -1
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = 1.875
$ws.Range("Y96").Value = -1
$ws.Range("Z96").Value = 0.825
$ws.Range("AA96").Value = -1
$ws.Range("AB96").Value = 0.9750000000000001

# Row 104
$ws.Range("B104").Value = 6537957
$ws.Range("E104").Value = "FC Flora Tallinn"
$ws.Range("F104").Value = "JK Nomme Kalju"
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = "D"
$ws.Range("J104").Value = 1.4
$ws.Range("K104").Value = 4
$ws.Range("L104").Value = 7.5
$ws.Range("M104").Value = 1.5
$ws.Range("N104").Value = 4.2
$ws.Range("O104").Value = 5
$ws.Range("P104").Value = -1
$ws.Range("Q104").Value = 1.85
$ws.Range("R104").Value = 1.95
$ws.Range("S104").Value = 2.75
$ws.Range("T104").Value = 1.85
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = -1
$ws.Range("W104").Value = 3.2
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.95
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.95

# Row 105
$ws.Range("B105").Value = 6533597
$ws.Range("E105").Value = "FC Kuressaare"
$ws.Range("F105").Value = "Parnu JK Vaprus"
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = "H"
$ws.Range("J105").Value = 2.5
$ws.Range("K105").Value = 3.4
$ws.Range("L105").Value = 2.5
$ws.Range("M105").Value = 2.15
$ws.Range("N105").Value = 3.6
$ws.Range("O105").Value = 2.875
$ws.Range("P105").Value = -0.25
$ws.Range("Q105").Value = 1.95
$ws.Range("R105").Value = 1.85
$ws.Range("S105").Value = 2.75
$ws.Range("T105").Value = 1.95
$ws.Range("U105").Value = 1.85
$ws.Range("V105").Value = 1.15
$ws.Range("W105").Value = -1
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = 0.95
$ws.Range("Z105").Value = -1
$ws.Range("AA105").Value = -1
$ws.Range("AB105").Value = 0.8500000000000001

# Row 106
$ws.Range("B106").Value = 6535416
$ws.Range("E106").Value = "Paide Linnameeskond"
$ws.Range("F106").Value = "FC Levadia Tallinn"
$ws.Range("G106").Value = 2
$ws.Range("H106").Value = 2
$ws.Range("I106").Value = "D"
$ws.Range("J106").Value = 3
$ws.Range("K106").Value = 3.8
$ws.Range("L106").Value = 2
$ws.Range("M106").Value = 3
$ws.Range("N106").Value = 4
$ws.Range("O106").Value = 1.909
$ws.Range("P106").Value = 0.5
$ws.Range("Q106").Value = 1.85
$ws.Range("R106").Value = 1.95
$ws.Range("S106").Value = 2.75
$ws.Range("T106").Value = 1.95
$ws.Range("U106").Value = 1.85
$ws.Range("V106").Value = -1
$ws.Range("W106").Value = 3
$ws.Range("X106").Value = -1
$ws.Range("Y106").Value = 0.8500000000000001
$ws.Range("Z106").Value = -1
$ws.Range("AA106").Value = 0.95
$ws.Range("AB106").Value = -1

# Row 107
$ws.Range("B107").Value = 6537869
$ws.Range("E107").Value = "JK Tallinna Kalev"
$ws.Range("F107").Value = "JK Trans Narva"
$ws.Range("G107").Value = 5
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = "H"
$ws.Range("J107").Value = 1.6
$ws.Range("K107").Value = 4
$ws.Range("L107").Value = 4.5
$ws.Range("M107").Value = 1.65
$ws.Range("N107").Value = 4
$ws.Range("O107").Value = 4.333
$ws.Range("P107").Value = -0.75
$ws.Range("Q107").Value = 1.8
$ws.Range("R107").Value = 2
$ws.Range("S107").Value = 2.75
$ws.Range("T107").Value = 1.9
$ws.Range("U107").Value = 1.9
$ws.Range("V107").Value = 0.6499999999999999
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = -1
$ws.Range("Y107").Value = 0.8
$ws.Range("Z107").Value = -1
$ws.Range("AA107").Value = 0.8999999999999999
$ws.Range("AB107").Value = -1

# Row 115
$ws.Range("B115").Value = 7919323
$ws.Range("E115").Value = "JK Nomme Kalju"
$ws.Range("F115").Value = "JK Trans Narva"
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = "H"
$ws.Range("J115").Value = 1.285
$ws.Range("K115").Value = 5.5
$ws.Range("L115").Value = 6.5
$ws.Range("M115").Value = 1.571
$ws.Range("N115").Value = 4.75
$ws.Range("O115").Value = 4.2
$ws.Range("P115").Value = -1
$ws.Range("Q115").Value = 1.925
$ws.Range("R115").Value = 1.875
$ws.Range("S115").Value = 2.75
$ws.Range("T115").Value = 1.875
$ws.Range("U115").Value = 1.925
$ws.Range("V115").Value = 0.571
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = -1
$ws.Range("Y115").Value = 0.925
$ws.Range("Z115").Value = -1
$ws.Range("AA115").Value = 0.4375
$ws.Range("AB115").Value = -0.5

# Row 116
$ws.Range("B116").Value = 7919322
$ws.Range("E116").Value = "FC Kuressaare"
$ws.Range("F116").Value = "FC Levadia Tallinn"
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 6
$ws.Range("I116").Value = "A"
$ws.Range("J116").Value = 11
$ws.Range("K116").Value = 6
$ws.Range("L116").Value = 1.166
$ws.Range("M116").Value = 15
$ws.Range("N116").Value = 8.5
$ws.Range("O116").Value = 1.125
$ws.Range("P116").Value = 2.5
$ws.Range("Q116").Value = 1.825
$ws.Range("R116").Value = 1.975
$ws.Range("S116").Value = 3.25
$ws.Range("T116").Value = 1.9
$ws.Range("U116").Value = 1.9
$ws.Range("V116").Value = -1
$ws.Range("W116").Value = -1
$ws.Range("X116").Value = 0.125
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = 0.9750000000000001
$ws.Range("AA116").Value = 0.8999999999999999
$ws.Range("AB116").Value = -1

# Row 120
$ws.Range("B120").Value = 7721007
$ws.Range("E120").Value = "JK Trans Narva"
$ws.Range("F120").Value = "JK Tammeka Tartu"
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 5
$ws.Range("I120").Value = "A"
$ws.Range("J120").Value = 2.25
$ws.Range("K120").Value = 3.3
$ws.Range("L120").Value = 2.75
$ws.Range("M120").Value = 2.1
$ws.Range("N120").Value = 3.25
$ws.Range("O120").Value = 3
$ws.Range("P120").Value = -0.25
$ws.Range("Q120").Value = 1.875
$ws.Range("R120").Value = 1.925
$ws.Range("S120").Value = 2.5
$ws.Range("T120").Value = 1.825
$ws.Range("U120").Value = 1.975
$ws.Range("V120").Value = -1
$ws.Range("W120").Value = -1
$ws.Range("X120").Value = 2
$ws.Range("Y120").Value = -1
$ws.Range("Z120").Value = 0.925
$ws.Range("AA120").Value = 0.825
$ws.Range("AB120").Value = -1

# Row 121
$ws.Range("B121").Value = 7721087
$ws.Range("E121").Value = "Paide Linnameeskond"
$ws.Range("F121").Value = "FC Flora Tallinn"
$ws.Range("G121").Value = 2
$ws.Range("H121").Value = 1
$ws.Range("I121").Value = "H"
$ws.Range("J121").Value = 2.2
$ws.Range("K121").Value = 3.3
$ws.Range("L121").Value = 2.8
$ws.Range("M121").Value = 1.85
$ws.Range("N121").Value = 3.6
$ws.Range("O121").Value = 3.4
$ws.Range("P121").Value = -0.5
$ws.Range("Q121").Value = 1.9
$ws.Range("R121").Value = 1.9
$ws.Range("S121").Value = 2.5
$ws.Range("T121").Value = 1.95
$ws.Range("U121").Value = 1.85
$ws.Range("V121").Value = 0.8500000000000001
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 0.8999999999999999
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.95
$ws.Range("AB121").Value = -1
